$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D/E/B/C column values in this sheet are text (coin price/volume
# strings and names), never real numbers -- force text format so Excel
# does not reinterpret numeric-looking strings as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.182.73"
$ws.Range("E2").Value = "  -5.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.232.10"
$ws.Range("E3").Value = "  -6.11%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.49"
$ws.Range("E5").Value = "  +0.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.21"
$ws.Range("E6").Value = "  -8.48%  "
$ws.Range("E7").Value = "  -8.36%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.563"
$ws.Range("E9").Value = "  -8.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "37.04"
$ws.Range("E10").Value = "  -10.19%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "54.24"
$ws.Range("E11").Value = "  -3.36%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0832"
$ws.Range("E12").Value = "  -9.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.68"
$ws.Range("E13").Value = "  -10.06%  "
$ws.Range("E14").Value = "  -1.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.868"
$ws.Range("E15").Value = "  -12.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.575.23"
$ws.Range("E16").Value = "  -6.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.44"
$ws.Range("E17").Value = "  -6.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.236.12"
$ws.Range("E18").Value = "  -6.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.116.66"
$ws.Range("E19").Value = "  -5.07%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.56"
$ws.Range("E20").Value = "  -5.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0₃0969"
$ws.Range("E21").Value = "  -9.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.56"
$ws.Range("E22").Value = "  -10.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "65.67"
$ws.Range("E23").Value = "  -10.78%  "
$ws.Range("E24").Value = "  -13.45%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "238.62"
$ws.Range("E25").Value = "  -10.70%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.17"
$ws.Range("E26").Value = "  -7.75%  "
$ws.Range("E27").Value = "  +0.21%  "
$ws.Range("E28").Value = "  +1.18%  "
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.03"
$ws.Range("E30").Value = "  -10.61%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.37"
$ws.Range("E31").Value = "  -15.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "35.99"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.40"
$ws.Range("E33").Value = "  -9.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0875"
$ws.Range("E34").Value = "  -7.97%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "153.85"
$ws.Range("E35").Value = "  -9.25%  "
$ws.Range("E36").Value = "  -5.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.17"
$ws.Range("E37").Value = "  +6.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.94"
$ws.Range("E38").Value = "  +1.05%  "
$ws.Range("E39").Value = "  -7.67%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.46"
$ws.Range("E40").Value = "  -5.87%  "
$ws.Range("E41").Value = "  -10.82%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.66"
$ws.Range("E42").Value = "  -8.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0325"
$ws.Range("E43").Value = "  -8.82%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.54"
$ws.Range("E44").Value = "  +4.38%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.736.84"
$ws.Range("E46").Value = "  -7.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "85.99"
$ws.Range("E47").Value = "  -12.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.205"
$ws.Range("E48").Value = "  -10.06%  "
$ws.Range("E49").Value = "  -10.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "75.53"
$ws.Range("E50").Value = "  -10.41%  "
$ws.Range("B51").Value = "FraxShare"
$ws.Range("C51").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.71"
$ws.Range("E51").Value = "  -5.90%  "
